# Add a new worksheet named "prompt" at the end of the workbook and
# populate its first cell with the prompt text, then make it the active
# (selected) sheet, matching the authored edit.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$promptSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$promptSheet.Name = "prompt"

$promptSheet.Range("A1").Value = "This is english prompt"

$promptSheet.Activate()
